# Update cryptos worksheet cell values per the commit diff (Mon Aug 26 19:36:24 UTC 2024 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.469.00'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.688.47'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'554.98"
$ws.Range("E5").Value = '  -3.70%  '
$ws.Range("D6").Value = "'158.37"
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -3.01%  '
$ws.Range("E9").Value = '  -4.12%  '
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("E11").Value = '  -4.49%  '
$ws.Range("D12").Value = "'5.37"
$ws.Range("E12").Value = '  -9.00%  '
$ws.Range("D13").Value = '3.162.73'
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").Value = "'26.37"
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").Value = '63.308.79'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("E16").Value = '  -4.05%  '
$ws.Range("D17").Value = '2.691.69'
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("D18").Value = "'12.04"
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").Value = "'4.57"
$ws.Range("E19").Value = '  -5.00%  '
$ws.Range("D20").Value = "'342.68"
$ws.Range("E20").Value = '  -4.62%  '
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = '  -4.93%  '
$ws.Range("D22").Value = "'0.995"
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  -4.07%  '
$ws.Range("D24").Value = "'63.92"
$ws.Range("E24").Value = '  -1.75%  '
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = "'8.16"
$ws.Range("E27").Value = '  -4.43%  '
$ws.Range("D28").Value = '0.0₃0855'
$ws.Range("E28").Value = '  -5.46%  '
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("E30").Value = '  -2.77%  '
$ws.Range("D31").Value = "'7.01"
$ws.Range("E31").Value = '  -4.86%  '
$ws.Range("D32").Value = "'165.21"
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("D34").Value = "'4.78"
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("D35").Value = "'19.53"
$ws.Range("E35").Value = '  -3.43%  '
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("D38").Value = "'340.04"
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("D39").Value = "'0.947"
$ws.Range("E39").Value = '  -5.89%  '
$ws.Range("D40").Value = "'6.06"
$ws.Range("E40").Value = '  -4.33%  '
$ws.Range("D41").Value = "'38.12"
$ws.Range("E41").Value = '  -2.47%  '
$ws.Range("E42").Value = '  -6.11%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'20.80"
$ws.Range("E43").Value = '  -5.46%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'20.29"
$ws.Range("E44").Value = '  -5.98%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = "'0.620"
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("D46").Value = "'0.0563"
$ws.Range("E46").Value = '  -4.34%  '
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = "'11.06"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = "'130.00"
$ws.Range("E49").Value = '  -5.36%  '
$ws.Range("E50").Value = '  -3.75%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0241"
$ws.Range("E51").Value = '  -5.04%  '
